$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Remember the currency number format used by the existing "Cost" column
$costFormat = $ws.Range("C2").NumberFormat

# Insert two new rows just above the current row 7 ("Transistors (x4)"),
# shifting the remaining rows (and their formulas) down by two rows.
$ws.Rows("7:8").Insert()

# Row 7: 3mm LED
$ws.Range("A7").Value = "3mm LED"
$ws.Range("B7").Value = "LED3R"
$ws.Range("C7").Value = 0.08
$ws.Range("C7").NumberFormat = $costFormat
$ws.Range("D7").Value = "Futurlec"
$ws.Range("E7").Value = "Ok"

# Row 8: 100R Resistor
$ws.Range("A8").Value = "100R Resistor"
$ws.Range("B8").Value = "R100R14W"
$ws.Range("C8").Formula = "=0.11/10"
$ws.Range("C8").NumberFormat = $costFormat
$ws.Range("D8").Value = "Futurlec"
$ws.Range("E8").Value = "Ok"

# Grow the table to include the two new rows (and the shifted-down totals row)
$lo.Resize($ws.Range("A1:E14"))
